$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values in columns B/C
$ws.Range("C2").Value = 12
$ws.Range("C3").Value = 11
$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.6
$ws.Range("C5").Value = 30

# Append a new (blank) row 6, duplicating row 5's formatting/style
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen column C slightly (content now fits a touch wider)
$ws.Columns.Item(3).ColumnWidth = 5.14
